$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: supplier changed from BuyaPi (40 Pin GPIO Connector Header, SKU: 412)
# to Creatron Inc. (7 Pin Receptacle Socket, CONHD-000007); price updated too.
# The hyperlink in H11 stays the same (same URL/text).
$ws.Range("A11").Value = "7 Pin Receptacle Socket"
$ws.Range("C11").Value = "CONHD-000007"
$ws.Range("B11").Value = "Creatron Inc."
$ws.Range("E11").Value = 0.45

# Row 12: Tax and Customs value updated
$ws.Range("E12").Value = 0.02

# Row 13: Shipping and Handling value now marked N/A instead of a number
$ws.Range("E13").Value = "N/A"

# Row 14: Grand total (CAD) recomputed
$ws.Range("G14").Value = 0.47

# Update the active selection to C14
$ws.Range("C14").Select()
